$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text values formatted like "265.69".
# Assigning a numeric-looking string directly would make Excel store it
# as a number, so we use a leading apostrophe to force text entry and
# then reset the style to Normal so no stray number-format/quote-prefix
# style sticks to the cell (matches the original plain inlineStr cells).
function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

# Simple price (column D) updates
Set-TextValue "D2" "265.69"
Set-TextValue "D3" "22.53"
Set-TextValue "D4" "6.271"
Set-TextValue "D5" "0.06158"
Set-TextValue "D7" "6.664"
Set-TextValue "D8" "1.341"
Set-TextValue "D9" "0.8290"
Set-TextValue "D11" "0.1588"
Set-TextValue "D12" "0.08297"
Set-TextValue "D14" "0.03182"
Set-TextValue "D40" "0.04620"
Set-TextValue "D41" "0.006956"
Set-TextValue "D42" "0.1139"
Set-TextValue "D43" "0.003130"
Set-TextValue "D44" "0.01080"
Set-TextValue "D45" "0.00006161"
Set-TextValue "D47" "0.7000"
Set-TextValue "D48" "0.1932"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.01240"

# Rows 15-26: coin reorder (rotation) with updated price and rank label
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09261"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.904"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001717"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04880"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006230"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005279"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001089"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.769"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.314"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D25" "0.3341"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D26" "0.1237"
$ws.Range("E26").Value = "25ProBitTokenPROB"
